$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet "pollutant" (sheet1): add 4 extra pollutant rows (split the
# PAH/PCB rows into separate short-name + long-name entries) and drop
# the trailing duplicate NMVOC/NMVOS row, giving a net +3 rows (37->40).
# ---------------------------------------------------------------------
$wsPollutant = $wb.Worksheets.Item("pollutant")
$loPollutant = $wsPollutant.ListObjects.Item(1)

# Grow the table by 3 rows so the table range + sheet can hold A1:B40.
$loPollutant.ListRows.Add() | Out-Null
$loPollutant.ListRows.Add() | Out-Null
$loPollutant.ListRows.Add() | Out-Null

# Final contents for rows 21..40 (rows 1..20 are unchanged).
$pollutantRows = @(
    @("PCDD-PCDF", "PCDD/ PCDF"),
    @("PCDD-PCDF", "PCDD/ PCDF`n(dioxins/ furans)"),
    @("BaP", "benzo(a) pyrene"),
    @("BbF", "benzo(b)"),
    @("BbF", "benzo(b) fluoranthene"),
    @("BkF", "benzo(k)"),
    @("BkF", "benzo(k) fluoranthene"),
    @("Indeno", "Indeno (1,2,3-cd) pyrene"),
    @("PAHs", "Total 1-4"),
    @("HCB", "HCB"),
    @("PCBs", "PCBs"),
    @("PCBs", "PCB"),
    @("NMVOC", "Totaal NMVOS"),
    @("NOx", "NOx (als NO2)"),
    @("NMVOC", "NMVOS"),
    @("NOx", "NOx"),
    @("PM2.5", "PM2,5"),
    @("SOx", "SOx (als SO2)"),
    @("SOx", "SOx"),
    @("PAHs", "PAK 4")
)

$r = 21
foreach ($pair in $pollutantRows) {
    $wsPollutant.Cells.Item($r, 1).Value = $pair[0]
    $wsPollutant.Cells.Item($r, 2).Value = $pair[1]
    $r = $r + 1
}

# Column widths widened to fit the longer pollutant names.
$wsPollutant.Columns(1).ColumnWidth = 22.75
$wsPollutant.Columns(2).ColumnWidth = 26.45

# ---------------------------------------------------------------------
# Sheet "nfr" (sheet2): rename the "prioriteit" table column to
# "NFR_priority".
# ---------------------------------------------------------------------
$wsNfr = $wb.Worksheets.Item("nfr")
$wsNfr.Range("C1").Value = "NFR_priority"
$wsNfr.Columns(3).ColumnWidth = 13.25
$wsNfr.Range("C1").Select()

# ---------------------------------------------------------------------
# Active sheet / selection: "pollutant" becomes the active tab (was
# "gnfr"), with A28 selected.
# ---------------------------------------------------------------------
$wsPollutant.Activate()
$wsPollutant.Range("A28").Select()
